$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, direct cell assignment matching the diff's before/after inline string values.

$ws.Range('D2').Value = '57.964.32'
$ws.Range('E2').Value = '  -1.44%  '

$ws.Range('D3').Value = '2.470.84'
$ws.Range('E3').Value = '  -1.04%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '517.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.54%  '

$ws.Range('E7').Value = '  +0.31%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.557'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.34%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0990'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.156'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.64%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.32%  '

$ws.Range('D13').Value = '2.912.41'
$ws.Range('E13').Value = '  -1.06%  '

$ws.Range('D14').Value = '57.901.15'
$ws.Range('E14').Value = '  -1.42%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.23%  '

$ws.Range('E16').Value = '  -2.25%  '

$ws.Range('D17').Value = '2.476.66'
$ws.Range('E17').Value = '  -1.31%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.72%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.28%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.24%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.23%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.25%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.46%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.409'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.82%  '

$ws.Range('E25').Value = '  +0.37%  '

$ws.Range('E26').Value = '  -2.71%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.95%  '

$ws.Range('D28').Value = '0.0₃0752'
$ws.Range('E28').Value = '  -2.58%  '

$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.37%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.33%  '

$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.78%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.72%  '

$ws.Range('E33').Value = '  +0.06%  '

$ws.Range('E34').Value = '  +0.44%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.16%  '

$ws.Range('E36').Value = '  -9.95%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.18%  '

$ws.Range('E38').Value = '  -4.57%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.788'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.26%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '272.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.60%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.592'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.42%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.75'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.73%  '

$ws.Range('E45').Value = '  -2.21%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0488'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.70%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0213'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.16%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '17.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.59%  '

$ws.Range('D49').Value = '1.734.48'
$ws.Range('E49').Value = '  -1.84%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.975'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.90%  '

$ws.Range('E51').Value = '  -0.94%  '
